$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K15").Value = -0.4017729932881683
$ws.Range("J16").Value = -0.04798648720847212
$ws.Range("I17").Value = -0.1961638776409175
$ws.Range("H18").Value = -0.3536098666450724
$ws.Range("G19").Value = 0.1470495898809984
$ws.Range("F20").Value = -0.1333082906544708
$ws.Range("E21").Value = 0.03866656719054083
$ws.Range("D22").Value = -0.2671959725557906
$ws.Range("C23").Value = 0.3451339801314955
$ws.Range("B24").Value = -0.343237405067616
